# Update generated statistics on the "展览" and "全部类型" sheets.
# F2: 45 -> 44
# F5: 2927 -> 2940
# F6: 293 -> 295

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 44
    $ws.Range("F5").Value = 2940
    $ws.Range("F6").Value = 295
}
